# New visitor attraction model
# Adds a second "random action interval / PB ramp sequence / ON-OFF /
# observation length" comparison block (columns L:Q) next to the existing
# reality-vs-simulation ratio table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new columns (L and O get an explicit width, like the authored file) ---
# ColumnWidth is expressed in characters; feed values tuned so the stored
# OOXML <col> width lands on (the closest attainable value to) 21.5 / 11.875.
$ws.Columns.Item(12).ColumnWidth = 20.714285714285715
$ws.Columns.Item(15).ColumnWidth = 11.142857142857142

# --- header row (row 1): L1 blank, M1/O1/Q1 labels ---
$ws.Range("M1").Value = "In reality"
$ws.Range("O1").Value = "In simulation"
$ws.Range("Q1").Value = "ratio"

# --- row 2: random action interval ---
$ws.Range("L2").Value = "random action interval "
$ws.Range("M2").Value = 2
$ws.Range("O2").Value = 0.5
$ws.Range("Q2").Formula = "=M2/O2"

# --- row 3: PB ramp sequence ---
$ws.Range("L3").Value = "PB ramp sequence"
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 15
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 3.75
$ws.Range("Q3").Formula = "=M3/O3"

# --- row 4: ON-OFF ---
$ws.Range("L4").Value = "ON-OFF "
$ws.Range("M4").Value = 1
$ws.Range("O4").Value = 0.25
$ws.Range("Q4").Formula = "=M4/O4"

# --- row 5: observation length ---
$ws.Range("L5").Value = "observation length "
$ws.Range("M5").Value = 4
$ws.Range("O5").Value = 1
$ws.Range("Q5").Formula = "=M5/O5"

# Vertically center the whole new block (adds the new cellXfs/alignment
# entry and stamps every cell in L1:Q5 with it, including the blanks).
$ws.Range("L1:Q5").VerticalAlignment = -4108

# --- view state: scroll so column G is left-most, select M13 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M13").Select()
